$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About": insert a new "Brazil GDP" source block (6 rows) right after
# the existing "Electricity, Heat, and CapEx Data" source block, pushing the
# "Notes" and "Amortized CapEx and OM Cost Notes" blocks further down.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("About")

$ws.Rows("14:19").Insert()

# B14: section header ("Brazil GDP") - reuse style of B5 (bold + grey fill)
$ws.Range("B14").Value = "Brazil GDP"
$ws.Range("B5").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B15: source name
$ws.Range("B15").Value = "Brazilian Institute of Geography and Statistics"

# B16: year - reuse style of B7 (number format)
$ws.Range("B16").Value = 2018
$ws.Range("B7").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B17: source detail
$ws.Range("B17").Value = "National System Account"

# B18: hyperlink to source
$ws.Range("B18").Value = "https://www.ibge.gov.br/estatisticas/economicas/contas-nacionais/9052-sistema-de-contas-nacionais-brasil.html?=&t=resultados"
$ws.Hyperlinks.Add($ws.Range("B18"), "https://www.ibge.gov.br/estatisticas/economicas/contas-nacionais/9052-sistema-de-contas-nacionais-brasil.html?=&t=resultados")
$ws.Range("B9").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet "Data": switch the country GDP reference from the U.S. to Brazil.
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$wsData.Range("A72").Value = "Brazil GDP"
$wsData.Range("B72").Value = 2.054
$wsData.Range("A74").Value = "Brazil GDP share"

# Downstream cells on Data / DACD-potential / DACD-energyintensity / DACD-capex
# sheets reference these values via formulas (B74, C78:H79, C83:H84, TREND(..))
# and will recalculate automatically.
